# Update "Ano" (date label, stored as text) and "Valor" (numeric) columns
# for the "Brasil" (rows 2-13) and "Sergipe" (rows 26-37) blocks.
# The "Nordeste" block (rows 14-25) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  Date = "01/12/2009"; Value = -13.26644370122631 },
    @{ Row = 3;  Date = "01/12/2010"; Value = -3.858520900321538 },
    @{ Row = 4;  Date = "01/12/2011"; Value = -6.700000000000006 },
    @{ Row = 5;  Date = "01/12/2012"; Value = -2.056807051909881 },
    @{ Row = 6;  Date = "01/12/2013"; Value = 3.654822335025365 },
    @{ Row = 7;  Date = "01/12/2014"; Value = 11.80476730987514 },
    @{ Row = 8;  Date = "01/12/2015"; Value = 7.177615571776141 },
    @{ Row = 9;  Date = "01/12/2016"; Value = -8.463251670378613 },
    @{ Row = 10; Date = "01/12/2017"; Value = -3.023758099352047 },
    @{ Row = 11; Date = "01/12/2018"; Value = -4.239917269906934 },
    @{ Row = 12; Date = "01/12/2019"; Value = -2.026342451874363 },
    @{ Row = 13; Date = "01/12/2020"; Value = 2.385892116182564 },
    @{ Row = 26; Date = "01/12/2009"; Value = -16.6147455867082 },
    @{ Row = 27; Date = "01/12/2010"; Value = 3.883495145631066 },
    @{ Row = 28; Date = "01/12/2011"; Value = -7.021063189568711 },
    @{ Row = 29; Date = "01/12/2012"; Value = -0.6972111553784854 },
    @{ Row = 30; Date = "01/12/2013"; Value = 1.006036217303818 },
    @{ Row = 31; Date = "01/12/2014"; Value = 24.87437185929651 },
    @{ Row = 32; Date = "01/12/2015"; Value = -2.331288343558291 },
    @{ Row = 33; Date = "01/12/2016"; Value = 0.8663366336633782 },
    @{ Row = 34; Date = "01/12/2017"; Value = -5.607476635514019 },
    @{ Row = 35; Date = "01/12/2018"; Value = 5.03067484662576 },
    @{ Row = 36; Date = "01/12/2019"; Value = -6.857142857142861 },
    @{ Row = 37; Date = "01/12/2020"; Value = 9.238451935081148 }
)

foreach ($u in $updates) {
    $dateCell = $ws.Cells.Item($u.Row, 3)
    # Force the "Ano" cell to stay plain text (it holds a DD/MM/YYYY label,
    # not a real date) instead of Excel auto-converting it to a date serial.
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $u.Date
    $dateCell.Style = "Normal"

    $ws.Cells.Item($u.Row, 4).Value = $u.Value
}
